$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6086
$ws.Range("K3").Value = 6273
$ws.Range("E4").Value = 2034
$ws.Range("K4").Value = 1312
$ws.Range("K6").Value = 6900
$ws.Range("E7").Value = 26040
$ws.Range("K7").Value = 21014

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 376
$ws.Range("K3").Value = 422
$ws.Range("K4").Value = 75
$ws.Range("K7").Value = 1380

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 464

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 331
$ws.Range("K6").Value = 277
$ws.Range("K7").Value = 914

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 204
$ws.Range("K3").Value = 238
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 716

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 166
$ws.Range("K7").Value = 496

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 90
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 56
$ws.Range("K7").Value = 612
$ws.Range("K8").Value = 1380
$ws.Range("K10").Value = 121
$ws.Range("K11").Value = 394
$ws.Range("K15").Value = 216
$ws.Range("K19").Value = 607
$ws.Range("K22").Value = 65
$ws.Range("K23").Value = 215
$ws.Range("K25").Value = 100
$ws.Range("K26").Value = 28
$ws.Range("K27").Value = 196
$ws.Range("K29").Value = 1143
$ws.Range("K31").Value = 235
$ws.Range("K33").Value = 914
$ws.Range("K36").Value = 272
$ws.Range("K37").Value = 716
$ws.Range("K42").Value = 779
$ws.Range("K43").Value = 176
$ws.Range("K45").Value = 28
$ws.Range("K48").Value = 264
$ws.Range("K51").Value = 271
$ws.Range("K52").Value = 553
$ws.Range("K54").Value = 411
$ws.Range("E55").Value = 246
$ws.Range("K59").Value = 36
$ws.Range("K60").Value = 124
$ws.Range("K65").Value = 496
$ws.Range("K67").Value = 821
$ws.Range("K72").Value = 103
$ws.Range("K76").Value = 285
$ws.Range("K78").Value = 237
$ws.Range("K79").Value = 525
$ws.Range("K82").Value = 23
$ws.Range("K83").Value = 464
$ws.Range("K85").Value = 973
$ws.Range("K89").Value = 307
$ws.Range("K91").Value = 240
$ws.Range("K92").Value = 81
$ws.Range("K93").Value = 79
$ws.Range("K94").Value = 284
$ws.Range("K99").Value = 345
$ws.Range("E101").Value = 26040
$ws.Range("K101").Value = 21014

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 235

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 295
$ws.Range("K6").Value = 235
$ws.Range("K7").Value = 821

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 224
$ws.Range("K7").Value = 411

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 330
$ws.Range("K3").Value = 412
$ws.Range("K6").Value = 318
$ws.Range("K7").Value = 1143

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 607

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 64
$ws.Range("K6").Value = 147
$ws.Range("K7").Value = 285

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 235
$ws.Range("K7").Value = 779

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("E4").Value = 22
$ws.Range("E7").Value = 246

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 62
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 115
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 175
$ws.Range("K3").Value = 172
$ws.Range("K7").Value = 525

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 106
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 202
$ws.Range("K6").Value = 165
$ws.Range("K7").Value = 612

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 102
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 394

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 86
$ws.Range("K7").Value = 307

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 41
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 317
$ws.Range("K4").Value = 55
$ws.Range("K7").Value = 973

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 23

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 150
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 553
